$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C4").Value = "783b3f4f-2bef-4f39-831d-f819f39bda8b"
$ws.Range("A11").Select()
